# Small update of input and buymodel dictionary for testing purposes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset Buy Rate values for the CPA vendor rows (27-31) to 0
$ws.Range("C27").Value2 = 0
$ws.Range("C28").Value2 = 0
$ws.Range("C29").Value2 = 0
$ws.Range("C30").Value2 = 0
$ws.Range("C31").Value2 = 0

# Add a new vendor row: ad4games / CPA / 0
$ws.Range("A32").Value2 = "ad4games"
$ws.Range("B32").Value2 = "CPA"
$ws.Range("C32").Value2 = 0

# Update the selection to match the new active cell used while editing
$ws.Range("F28").Select() | Out-Null
